$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Purchase 22-23")
$ws.Range("A1").Value = "test"
